$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws2 = $wb.Worksheets.Item("ProductLoanOutput")

# Update the product name text (adds a hyphen after "198") on both sheets
$ws1.Range("B1").Value = "198-MS-EI-DB-DL-REC-NON-RNI-CTPD-DL-MD-TR-1-ONTIME"
$ws2.Range("B1").Value = "198-MS-EI-DB-DL-REC-NON-RNI-CTPD-DL-MD-TR-1-ONTIME"

# Change selection on the input sheet
$ws1.Range("B1").Select()

# Activate the output sheet and set its selection, making it the active tab
$ws2.Activate()
$ws2.Range("B1").Select()
